$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Monthly refresh: the report now covers one more month (November 2016).
# This:
#   1) inserts a new "November" data row into the "Year 2016" block (after
#      October, before the "Year to Date" summary block),
#   2) rolls the "Year to Date" (2014/2015/2016) totals forward one month,
#   3) rolls the "Rolling 12 Months Ending in ..." block forward one month
#      (relabeling it, and updating its 2015/2016 totals), and
#   4) updates the report's subtitle to say "November 2016" instead of
#      "October 2016".
# ---------------------------------------------------------------------------

# 1) Insert a new row above the old row 52 ("Year to Date" header), which
#    pushes everything at/after row 52 down by one row.
$ws.Rows("52").Insert()

# Give the new row the same look as the other month rows (copy the format
# of the "October" row right above it: label style + right-aligned number
# style for the data columns).
$ws.Range("A51:F51").Copy()
$ws.Range("A52:F52").PasteSpecial(-4122)

# Fill in the November 2016 figures.
$ws.Range("A52").Value = "November"
$ws.Range("B52").Value = 11829
$ws.Range("C52").Value = 10707
$ws.Range("D52").Value = 4985
$ws.Range("E52").Value = 53
$ws.Range("F52").Value = 27574

# 2) "Year to Date" totals (now rows 54-56) updated through November 2016.
$ws.Range("B54").Value = 161471
$ws.Range("C54").Value = 134061
$ws.Range("D54").Value = 65483
$ws.Range("E54").Value = 745
$ws.Range("F54").Value = 361760

$ws.Range("B55").Value = 163864
$ws.Range("C55").Value = 133956
$ws.Range("D55").Value = 63123
$ws.Range("E55").Value = 710
$ws.Range("F55").Value = 361654

$ws.Range("B56").Value = 161836
$ws.Range("C56").Value = 129903
$ws.Range("D56").Value = 58213
$ws.Range("E56").Value = 649
$ws.Range("F56").Value = 350601

# 3) "Rolling 12 Months Ending in ..." block (now row 57 header, rows 58-59
#    data) relabeled and rolled forward to November.
$ws.Range("A57").Value = "Rolling 12 Months Ending in November"

$ws.Range("B58").Value = 178571
$ws.Range("C58").Value = 145148
$ws.Range("D58").Value = 68496
$ws.Range("E58").Value = 775
$ws.Range("F58").Value = 392990

$ws.Range("B59").Value = 175595
$ws.Range("C59").Value = 140728
$ws.Range("D59").Value = 63256
$ws.Range("E59").Value = 710
$ws.Range("F59").Value = 380288

# 4) Subtitle: "... 2006 - October 2016 ..." -> "... 2006 - November 2016 ...".
$ws.Range("A2").Value = "Total by End-Use Sector, 2006 - November 2016 (Million Dollars)"
